$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Config -------------------------------------------------------------
$firstRow = 836
$lastRow  = 942          # last row of the existing block (before the edit)
$newLastRow = $lastRow + 2  # two brand-new rows get appended at the end

# Columns that vary row-to-row within this block.
# (A,B,C,E,F,G,H,N,Q,R are identical for every row in the block.)
$varCols = @(4, 9, 10, 11, 12, 13, 15, 16)   # D, I, J, K, L, M, O, P

# --- 1. Snapshot every row of the existing block ------------------------
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @{}
    foreach ($c in $varCols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# Full row snapshot (all columns A..R) for the two rows that will be
# pushed past the old last row (they need every column written out,
# since those target rows currently have no content at all).
$fullCols = 1..18
$fullSnapshot = @{}
foreach ($r in @($lastRow - 1, $lastRow)) {
    $rowData = @{}
    foreach ($c in $fullCols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $fullSnapshot[$r] = $rowData
}

# --- 2. Write the two brand-new records into rows 836 & 837 -------------
# New record 1: fecha 45142, Primera
$ws.Cells.Item(836, 4).Value  = 45142
$ws.Cells.Item(836, 9).Value  = "Primera"
$ws.Cells.Item(836, 10).Value = 70
$ws.Cells.Item(836, 11).Value = 12000
$ws.Cells.Item(836, 12).Value = 12000
$ws.Cells.Item(836, 13).Value = 12000
$ws.Cells.Item(836, 15).Value = "Región Metropolitana"
$ws.Cells.Item(836, 16).Value = 4000

# New record 2: fecha 45142, Segunda
$ws.Cells.Item(837, 4).Value  = 45142
$ws.Cells.Item(837, 9).Value  = "Segunda"
$ws.Cells.Item(837, 10).Value = 52
$ws.Cells.Item(837, 11).Value = 9000
$ws.Cells.Item(837, 12).Value = 9000
$ws.Cells.Item(837, 13).Value = 9000
$ws.Cells.Item(837, 15).Value = "Región Metropolitana"
$ws.Cells.Item(837, 16).Value = 3000

# --- 3. Shift the rest of the block down by two rows ---------------------
# new row (r) = old row (r - 2), for r = 838 .. 942
for ($r = ($firstRow + 2); $r -le $lastRow; $r++) {
    $src = $snapshot[$r - 2]
    foreach ($c in $varCols) {
        $ws.Cells.Item($r, $c).Value = $src[$c]
    }
}

# --- 4. Append the final two (previously last) rows at 943 & 944 --------
$destRows = @($newLastRow - 1, $newLastRow)   # 943, 944
$srcRows  = @($lastRow - 1, $lastRow)         # 941, 942

for ($i = 0; $i -lt $destRows.Length; $i++) {
    $dest = $destRows[$i]
    $src  = $srcRows[$i]
    $data = $fullSnapshot[$src]
    foreach ($c in $fullCols) {
        $ws.Cells.Item($dest, $c).Value = $data[$c]
    }
    # Make sure the date column keeps the date number format.
    $ws.Cells.Item($dest, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Host "done"
